$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.953.44"
Set-TextValue $ws.Range("E2") "  +5.39%  "
Set-TextValue $ws.Range("D3") "3.593.09"
Set-TextValue $ws.Range("E3") "  +5.17%  "
Set-TextValue $ws.Range("E4") "  -0.01%  "
Set-TextValue $ws.Range("D5") "587.19"
Set-TextValue $ws.Range("E5") "  +3.34%  "
Set-TextValue $ws.Range("D6") "190.00"
Set-TextValue $ws.Range("E6") "  +4.52%  "
Set-TextValue $ws.Range("D7") "0.644"
Set-TextValue $ws.Range("E7") "  +1.82%  "
Set-TextValue $ws.Range("D8") "3.587.23"
Set-TextValue $ws.Range("E8") "  +5.15%  "
Set-TextValue $ws.Range("E9") "  -0.07%  "
Set-TextValue $ws.Range("D10") "0.177"
Set-TextValue $ws.Range("E10") "  -0.36%  "
Set-TextValue $ws.Range("D11") "0.659"
Set-TextValue $ws.Range("E11") "  +2.87%  "
Set-TextValue $ws.Range("D12") "57.93"
Set-TextValue $ws.Range("E12") "  +5.75%  "
Set-TextValue $ws.Range("D13") "0.0000290"
Set-TextValue $ws.Range("E13") "  +3.35%  "
Set-TextValue $ws.Range("D14") "9.74"
Set-TextValue $ws.Range("E14") "  +4.27%  "
Set-TextValue $ws.Range("D15") "4.164.74"
Set-TextValue $ws.Range("E15") "  +5.03%  "
Set-TextValue $ws.Range("B16") "WrappedEther"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D16") "3.605.52"
Set-TextValue $ws.Range("E16") "  +5.54%  "
Set-TextValue $ws.Range("B17") "Chainlink"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D17") "19.29"
Set-TextValue $ws.Range("E17") "  +4.97%  "
Set-TextValue $ws.Range("D18") "69.888.49"
Set-TextValue $ws.Range("E18") "  +5.45%  "
Set-TextValue $ws.Range("D19") "12.43"
Set-TextValue $ws.Range("E19") "  +3.50%  "
Set-TextValue $ws.Range("D20") "0.120"
Set-TextValue $ws.Range("E20") "  +0.15%  "
Set-TextValue $ws.Range("E21") "  +4.08%  "
Set-TextValue $ws.Range("D22") "493.74"
Set-TextValue $ws.Range("E22") "  +4.80%  "
Set-TextValue $ws.Range("D23") "17.54"
Set-TextValue $ws.Range("E23") "  +19.45%  "
Set-TextValue $ws.Range("D24") "5.38"
Set-TextValue $ws.Range("E24") "  +7.50%  "
Set-TextValue $ws.Range("D25") "4.44"
Set-TextValue $ws.Range("E25") "  +7.18%  "
Set-TextValue $ws.Range("D26") "90.52"
Set-TextValue $ws.Range("E26") "  +1.08%  "
Set-TextValue $ws.Range("D27") "3.09"
Set-TextValue $ws.Range("E27") "  +5.00%  "
Set-TextValue $ws.Range("D28") "11.07"
Set-TextValue $ws.Range("E28") "  +1.88%  "
Set-TextValue $ws.Range("D29") "9.39"
Set-TextValue $ws.Range("E29") "  +5.84%  "
Set-TextValue $ws.Range("D30") "32.19"
Set-TextValue $ws.Range("E30") "  +2.51%  "
Set-TextValue $ws.Range("D31") "7.66"
Set-TextValue $ws.Range("E31") "  +10.21%  "
Set-TextValue $ws.Range("D32") "12.20"
Set-TextValue $ws.Range("E32") "  +5.32%  "
Set-TextValue $ws.Range("D33") "619.35"
Set-TextValue $ws.Range("E33") "  +5.63%  "
Set-TextValue $ws.Range("E34") "  +7.03%  "
Set-TextValue $ws.Range("D35") "65.10"
Set-TextValue $ws.Range("E35") "  +4.11%  "
Set-TextValue $ws.Range("D36") "0.0₃0819"
Set-TextValue $ws.Range("E36") "  +7.96%  "
Set-TextValue $ws.Range("D37") "38.08"
Set-TextValue $ws.Range("E37") "  +4.44%  "
Set-TextValue $ws.Range("B38") "Dai"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D38") "1.00"
Set-TextValue $ws.Range("E38") "  +0.05%  "
Set-TextValue $ws.Range("B39") "TheGraph"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D39") "0.403"
Set-TextValue $ws.Range("E39") "  +4.46%  "
Set-TextValue $ws.Range("B40") "Kaspa"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D40") "0.146"
Set-TextValue $ws.Range("E40") "  -0.23%  "
Set-TextValue $ws.Range("E41") "  +0.60%  "
Set-TextValue $ws.Range("D42") "3.313.74"
Set-TextValue $ws.Range("E42") "  +5.78%  "
Set-TextValue $ws.Range("D43") "3.07"
Set-TextValue $ws.Range("E43") "  +5.46%  "
Set-TextValue $ws.Range("D44") "0.0445"
Set-TextValue $ws.Range("E44") "  +4.72%  "
Set-TextValue $ws.Range("D45") "2.67"
Set-TextValue $ws.Range("E45") "  +5.07%  "
Set-TextValue $ws.Range("D46") "3.27"
Set-TextValue $ws.Range("E46") "  +2.26%  "
Set-TextValue $ws.Range("D47") "0.137"
Set-TextValue $ws.Range("E47") "  +1.85%  "
Set-TextValue $ws.Range("D48") "9.11"
Set-TextValue $ws.Range("E48") "  +5.93%  "
Set-TextValue $ws.Range("B49") "LidoDAOToken"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D49") "3.32"
Set-TextValue $ws.Range("E49") "  +5.18%  "
Set-TextValue $ws.Range("B50") "dogwifhat"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D50") "2.69"
Set-TextValue $ws.Range("E50") "  -4.00%  "
Set-TextValue $ws.Range("D51") "1.00"
Set-TextValue $ws.Range("E51") "  +0.09%  "

Write-Host "Applied 109 cell updates"
